# Loan RBI, Variable Instalments
#
# The "Repayment schedule" sheet gains a new (blank) column between the
# existing "Late" column (N) and the "Outstanding" column (which, after the
# insert, moves from N to O). All of the old N/O/P columns (Late, Outstanding,
# heading/Date) shift one column to the right (-> O/P/Q), and the sheet
# becomes the active tab/selection of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Make this sheet the active one (workbook.activeTab + sheetView.tabSelected
# both follow from activating it).
$ws.Activate()

# Insert a new blank column at N, pushing the old N/O/P columns to O/P/Q.
# Excel normally carries over the column width/format of the column to the
# left (M) onto the freshly inserted column, so mirror that explicitly.
$mColumnWidth = $ws.Columns("M").ColumnWidth
$ws.Columns("N").Insert()
$ws.Columns("N").ColumnWidth = $mColumnWidth

# Update the selected cell on this sheet to match the new layout.
$ws.Range("K15").Select()
